# aula 07 - INDMO
# Fill in the attendance ("P" = presente, "F" = falta) for the two new
# date columns (BC = 2022-12-05, BD = 2022-12-06) and tidy up the header
# row / selection / column defaults to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chamada")

$xlCenter = -4108
$xlLeft   = -4131

# ---------------------------------------------------------------------
# 1) Header row (row 1) date-cell formatting: BA1 and BD1:BM1 move from
#    the plain / "special" date styles onto the common centered-date
#    style already used by the rest of the row (AJ1:AZ1).
# ---------------------------------------------------------------------
$ws.Range("BA1").HorizontalAlignment = $xlCenter
$ws.Range("BD1:BM1").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------
# 3) Attendance data for class "aula 07" (INDMO), columns BC (05/12) and
#    BD (06/12), rows 3-30 (row 13 is a hidden/unused student row).
# ---------------------------------------------------------------------
$attendance = @{
    3  = @("P","P")
    4  = @("P","P")
    5  = @("P","F")
    6  = @("P","P")
    7  = @("P","P")
    8  = @("P","P")
    9  = @("P","P")
    10 = @("P","P")
    11 = @("P","P")
    12 = @("F","P")
    14 = @("P","P")
    15 = @("P","P")
    16 = @("P","F")
    17 = @("P","P")
    18 = @("P","P")
    19 = @("P","P")
    20 = @("P","P")
    21 = @("P","P")
    22 = @("P","P")
    23 = @("P","F")
    24 = @("P","P")
    25 = @("P","P")
    26 = @("P","P")
    27 = @("P","F")
    28 = @("F","P")
    29 = @("P","F")
    30 = @("P","P")
}

foreach ($row in $attendance.Keys) {
    $vals = $attendance[$row]
    $bc = $ws.Cells.Item($row, 55)   # column BC
    $bd = $ws.Cells.Item($row, 56)   # column BD

    $bc.Value = $vals[0]

    $bd.Value = $vals[1]
    $bd.HorizontalAlignment = $xlLeft
}

# ---------------------------------------------------------------------
# 4) Restore the selection Excel leaves behind after entering the last
#    bit of attendance data.
# ---------------------------------------------------------------------
$ws.Range("BD29").Select()
